$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 581
$ws.Range("F10").Value = 286
$ws.Range("H10").Value = 382

$ws.Range("E11").Value = 366
$ws.Range("F11").Value = 196
$ws.Range("H11").Value = 260

$ws.Range("E12").Value = 565
$ws.Range("F12").Value = 301
$ws.Range("H12").Value = 387

$ws.Range("E13").Value = 140

$ws.Range("E14").Value = 126

$ws.Range("E15").Value = 173

$ws.Range("E17").Value = 101

$ws.Range("E19").Value = 15

$ws.Range("E24").Value = 218
$ws.Range("F24").Value = 115
$ws.Range("H24").Value = 145

$ws.Range("E25").Value = 276
$ws.Range("F25").Value = 138
$ws.Range("H25").Value = 198

$ws.Range("E27").Value = 335
$ws.Range("F27").Value = 171
$ws.Range("H27").Value = 252

$ws.Range("E28").Value = 202

$ws.Range("E30").Value = 215

$ws.Range("E33").Value = 300
$ws.Range("F33").Value = 154
$ws.Range("H33").Value = 243

$ws.Range("E34").Value = 222
$ws.Range("F34").Value = 147
$ws.Range("H34").Value = 186

$ws.Range("E38").Value = 93

$ws.Range("E39").Value = 181

$ws.Range("E40").Value = 267
$ws.Range("F40").Value = 125
$ws.Range("H40").Value = 205

$ws.Range("E41").Value = 395

$ws.Range("E42").Value = 387
$ws.Range("F42").Value = 216
$ws.Range("H42").Value = 277

$ws.Range("E44").Value = 313

$ws.Range("E45").Value = 148
$ws.Range("F45").Value = 73
$ws.Range("H45").Value = 112

$ws.Range("E46").Value = 329

$ws.Range("E47").Value = 462
$ws.Range("F47").Value = 237
$ws.Range("H47").Value = 329

$ws.Range("E49").Value = 289

$ws.Range("E50").Value = 244

$ws.Range("F51").Value = 108
$ws.Range("H51").Value = 182
